# The workbook was re-saved (selection state updated) after being opened on
# another device/build of Excel (per commit message, on iPad). The only
# user-visible/model-visible change is the current selection on the first
# worksheet: entire columns H:W are now selected with H1 as the active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Select() | Out-Null
$ws.Range("H1:W1048576").Select() | Out-Null
